$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2o Parcial")

# Row 8
$ws.Range("E8").Value = 20
$ws.Range("F8").Value = 4
$ws.Range("G8").Value = 83.3
$ws.Range("H8").Value = 16.7
$ws.Range("I8").Value = 7.6
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0

# Row 13
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0

# Row 14
$ws.Range("E14").Value = 182
$ws.Range("F14").Value = 9
$ws.Range("G14").Value = 95.3
$ws.Range("H14").Value = 4.7
$ws.Range("I14").Value = 8.699999999999999
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0

# Row 15
$ws.Range("E15").Value = 260
$ws.Range("F15").Value = 12
$ws.Range("G15").Value = 95.59999999999999
$ws.Range("H15").Value = 4.4
$ws.Range("I15").Value = 8.4
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
